# IntroToNodeJS.pptx edit script
# - Slide 37 ("Install and update Node packages..." title slide): drop the
#   stray trailing empty-paragraph mark (<a:endParaRPr/>) after the title run.
# - Slide 5 ("About Node"): merge the "Windows, Linux, Mac OSX" bullet and the
#   "Still in "beta" phase" bullet into a single paragraph, splitting the
#   surviving text into three runs ("Windows, Linux, " / "Mac " / "OSX") and
#   dropping the old trailing empty paragraph.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 37: "Install and update Node packages through the Node Package Manager"
# ---------------------------------------------------------------------------
$slide37 = $p.Slides.Item(37)
$titleShape = $slide37.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# Clearing the range fully and re-typing the same text drops the leftover
# paragraph-end run properties that otherwise linger on the paragraph mark.
$titleText = $titleRange.Text
$titleRange.Delete()
$titleRange.Text = $titleText

# ---------------------------------------------------------------------------
# Slide 5: "About Node" bullet list
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$contentShape = $slide5.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange

# Paragraph 5 = "Windows, Linux, Mac OSX"
# Paragraph 6 = "Still in "beta" phase"      <- remove entirely
# Paragraph 7 = empty trailing paragraph     <- remove entirely, merging its
#                                                paragraph mark into para 5
$betaParagraph = $contentRange.Paragraphs(6)
$betaParagraph.Delete()

$winParagraph = $contentRange.Paragraphs(5)
$firstRun = $winParagraph.Runs(1)
$firstRun.Text = "Windows, Linux, "
$winParagraph.InsertAfter("Mac ") | Out-Null
$winParagraph.InsertAfter("OSX") | Out-Null

$trailingParagraph = $contentRange.Paragraphs(6)
$trailingParagraph.Delete()
